$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Form-feed character used throughout the original invoice text runs
# (serialized back out to OOXML as the literal escape "_x000C_").
$ff = [char]0x0C

# Drop the now-unused last column (G) and the two extra invoice-header /
# invoice-footer rows (4 and 5); row 1 is left untouched so rows 2/3 keep
# their original row numbers.
$ws.Columns("G").Delete()
$ws.Rows(5).Delete()
$ws.Rows(4).Delete()

# Row 2 becomes the line-item column headers, row 3 the corresponding
# values - content that used to live packed into the old (now-deleted)
# row 4 cells, split one-value-per-cell. Filled column-by-column (header
# then value) to mirror how the line-item block reads top-to-bottom.
$ws.Range("B2").Value = " Hrs/Qty`n" + $ff

# " 1.00" alone round-trips through plain cell types as a *number* (it is
# indistinguishable from the numeric literal 1), but the source workbook
# stores it as text. Route it through a formula + paste-values so Excel
# keeps the string type without disturbing the cell's existing style.
$ws.Range("B3").Formula = '=" 1.00"&CHAR(10)&CHAR(' + [int]$ff + ')'
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("C2").Value = " Service`n" + $ff
$ws.Range("C3").Value = " Web Design`nThis is a sample description...`n" + $ff

$ws.Range("D2").Value = " Rate/Price`n" + $ff
$ws.Range("D3").Value = " `$85.00`n" + $ff

$ws.Range("E2").Value = " Adjust`n" + $ff
$ws.Range("E3").Value = " 0.00%`n" + $ff

$ws.Range("F2").Value = " Sub Total`n" + $ff
$ws.Range("F3").Value = " `$85.00`n" + $ff

# Multi-line text triggers automatic row-height adjustment; auto-fitting
# restores the rows to their default (no explicit custom height), which
# is what the target worksheet looks like.
$ws.Rows(2).AutoFit()
$ws.Rows(3).AutoFit()
